# "Generate Report for Archive"
#
# 1) The localization status string changes from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F4 and the
#    Status column - column C - on the zh-cn and de-de report sheets).
# 2) The two "Status" columns on the Overview sheet (zh-cn / de-de,
#    columns E & F) and the "Status" column (column C) on the zh-cn and
#    de-de detail sheets are narrowed now that the status text is
#    shorter ("In Translation" instead of "Ready for handoff").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------

foreach ($cellRef in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    $cell = $overview.Range($cellRef)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value2 = $newStatus
    }
}

foreach ($ws in @($zhcn, $dede)) {
    foreach ($cellRef in @("C2", "C3", "C4")) {
        $cell = $ws.Range($cellRef)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}

# --- Narrow the Status columns to fit the shorter text ----------------------
# (ColumnWidth is expressed in characters of the Normal style font; Excel
# snaps the stored column width to whole pixels, so this lands the column on
# the narrowest pixel width available once "In Translation" replaces the
# longer "Ready for handoff" text.)

$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # Overview!E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # Overview!F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # de-de!C (Status)
